# Applies the cryptos-list price/volume update described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '27.093.54'
$ws.Range('E2').Value = '  -2.59%  '

# Row 3
$ws.Range('D3').Value = '1.867.94'
$ws.Range('E3').Value = '  -2.02%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9998'
$ws.Range('E4').Value = '  -0.16%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '306.18'
$ws.Range('E5').Value = '  -2.16%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9996'

# Row 7
$ws.Range('E7').Value = '  -1.23%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3770'
$ws.Range('E8').Value = '  -0.42%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07164'
$ws.Range('E9').Value = '  -1.02%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8907'
$ws.Range('E10').Value = '  -2.33%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.69'
$ws.Range('E11').Value = '  -2.75%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07605'
$ws.Range('E12').Value = '  -0.53%  '

# Row 13
$ws.Range('D13').Value = '1.880.60'
$ws.Range('E13').Value = '  -1.64%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.313'
$ws.Range('E14').Value = '  -2.62%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '89.75'
$ws.Range('E15').Value = '  -2.61%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.000'
$ws.Range('E16').Value = '  -0.19%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008494'
$ws.Range('E17').Value = '  -2.39%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '14.08'
$ws.Range('E18').Value = '  -3.19%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.9998'
$ws.Range('E19').Value = '  -0.11%  '

# Row 20
$ws.Range('D20').Value = '27.113.44'
$ws.Range('E20').Value = '  -2.63%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.030'
$ws.Range('E21').Value = '  -2.43%  '

# Row 22
$ws.Range('D22').Value = '2.117.11'
$ws.Range('E22').Value = '  -2.11%  '

# Row 23
$ws.Range('E23').Value = '  -3.23%  '

# Row 24
$ws.Range('E24').Value = '  -2.70%  '

# Row 25
$ws.Range('E25').Value = '  -1.80%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '147.51'
$ws.Range('E26').Value = '  -4.16%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.97'
$ws.Range('E27').Value = '  -1.95%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.096'
$ws.Range('E28').Value = '  -3.39%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '112.84'
$ws.Range('E29').Value = '  -1.74%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.665'
$ws.Range('E30').Value = '  -4.13%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.677'
$ws.Range('E31').Value = '  -3.91%  '

# Row 32
$ws.Range('E32').Value = '  +1.46%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05113'
$ws.Range('E33').Value = '  -3.25%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.070'
$ws.Range('E34').Value = '  -3.42%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.158'
$ws.Range('E35').Value = '  -6.61%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7263'
$ws.Range('E36').Value = '  -6.97%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02034'
$ws.Range('E37').Value = '  -2.85%  '

# Row 38
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.506'
$ws.Range('E38').Value = '  -3.94%  '

# Row 39
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.062'
$ws.Range('E39').Value = '  -0.43%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.073'
$ws.Range('E40').Value = '  -1.87%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.5326'
$ws.Range('E41').Value = '  -4.67%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.495'
$ws.Range('E42').Value = '  -3.53%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '116.40'
$ws.Range('E43').Value = '  +0.52%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.296'
$ws.Range('E44').Value = '  -3.03%  '

# Row 45
$ws.Range('E45').Value = '  -3.21%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4630'
$ws.Range('E46').Value = '  -3.77%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.9990'
$ws.Range('E47').Value = '  -0.09%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.979'
$ws.Range('E48').Value = '  -5.03%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.573'
$ws.Range('E49').Value = '  -3.12%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '36.54'
$ws.Range('E50').Value = '  -1.41%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '63.66'
$ws.Range('E51').Value = '  -5.04%  '
